# Quarterly balance sheet update: roll the 5 visible quarters forward by one
# (drop the oldest quarter from column D, shift D<-E<-F<-G<-H, and populate the
# new column H with the newest quarter's figures), refresh the period / publish
# date headers to match, bump the copyright year, and nudge the column widths
# that travelled along with the shifted data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copyright year bump -------------------------------------------------
$ws.Range("B3").Value = "Copyright @2015 - 2023"

# --- Column headers (row 8: فصل / quarter-end labels) --------------------
$ws.Range("D8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1401/09"

# --- Column headers (row 9: تاریخ انتشار / publish date labels) ----------
$ws.Range("D9").Value = "1400-11-02"
$ws.Range("E9").Value = "1401-11-01 (8)"
$ws.Range("F9").Value = "1401-05-01"
$ws.Range("G9").Value = "1401-08-30 (2)"
$ws.Range("H9").Value = "1401-11-01"

# --- Column widths (D keeps 29, E<-31, F/G/H settle at 29) ---------------
$ws.Columns.Item(5).ColumnWidth = 30.17
$ws.Columns.Item(6).ColumnWidth = 28.17

# --- Data rows: shift D<-E<-F<-G<-H, H <- new quarter value --------------
$rowNewH = @(
    @{Row=12; NewH=1708760}
    @{Row=13; NewH=3758124}
    @{Row=14; NewH=400577}
    @{Row=15; NewH=4348400}
    @{Row=16; NewH=630328}
    @{Row=17; NewH=$null}
    @{Row=18; NewH=10846189}
    @{Row=19; NewH=$null}
    @{Row=20; NewH=922755}
    @{Row=21; NewH=$null}
    @{Row=22; NewH=633696}
    @{Row=23; NewH=60391}
    @{Row=24; NewH=$null}
    @{Row=25; NewH=$null}
    @{Row=26; NewH=1616842}
    @{Row=27; NewH=12463031}
    @{Row=29; NewH=925703}
    @{Row=30; NewH=$null}
    @{Row=31; NewH=714195}
    @{Row=32; NewH=1655015}
    @{Row=33; NewH=179836}
    @{Row=34; NewH=17160}
    @{Row=35; NewH=$null}
    @{Row=36; NewH=$null}
    @{Row=37; NewH=3491909}
    @{Row=38; NewH=$null}
    @{Row=39; NewH=$null}
    @{Row=40; NewH=0}
    @{Row=41; NewH=554566}
    @{Row=42; NewH=554566}
    @{Row=43; NewH=4046475}
    @{Row=45; NewH=$null}
    @{Row=46; NewH=$null}
    @{Row=47; NewH=$null}
    @{Row=48; NewH=-13}
    @{Row=49; NewH=8458}
    @{Row=50; NewH=$null}
    @{Row=51; NewH=$null}
    @{Row=52; NewH=$null}
    @{Row=53; NewH=$null}
    @{Row=54; NewH=$null}
    @{Row=55; NewH=$null}
    @{Row=56; NewH=7693111}
    @{Row=57; NewH=8416556}
    @{Row=58; NewH=12463031}
)

foreach ($entry in $rowNewH) {
    $r = $entry.Row
    $rng = $ws.Range("D" + $r + ":H" + $r)
    $vals = $rng.Value2

    $d = $vals[1,1]
    $e = $vals[1,2]
    $f = $vals[1,3]
    $g = $vals[1,4]
    $h = $vals[1,5]

    if ($null -ne $entry.NewH) {
        $newH = $entry.NewH
    } else {
        $newH = $h
    }

    $out = New-Object 'object[,]' 1,5
    $out[0,0] = $e
    $out[0,1] = $f
    $out[0,2] = $g
    $out[0,3] = $h
    $out[0,4] = $newH

    $rng.Value = $out
}
